# Apply the edit described by the diff:
#  1. Rename worksheet from "GammaFiber2F-HW35.xpc" to "GammaFiber2F"
#  2. Append a new data row (row 16) with HKL index 14 / label "HexGrid-60degTilt5degRes"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet
$ws.Name = "GammaFiber2F"

# 2. Append the new row of averaged intensity data
$row = 16

$ws.Cells.Item($row, 1).Value = 14
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Cells.Item($row, 2).Value = "HexGrid-60degTilt5degRes"

$ws.Cells.Item($row, 3).Value  = 1.017444091536401
$ws.Cells.Item($row, 4).Value  = 0.9216615679867284
$ws.Cells.Item($row, 5).Value  = 1.008359098580538
$ws.Cells.Item($row, 6).Value  = 1.017444091536401
$ws.Cells.Item($row, 7).Value  = 0.9606089018613844
$ws.Cells.Item($row, 8).Value  = 1.041424079514018
$ws.Cells.Item($row, 9).Value  = 1.010712039757009
$ws.Cells.Item($row, 10).Value = 0.9216615679867284
$ws.Cells.Item($row, 11).Value = 0.9650103332836334
$ws.Cells.Item($row, 12).Value = 0.9912272124100172
$ws.Cells.Item($row, 13).Value = 0.9933682965393466
